# -----------------------------------------------------------------------
# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") has a comparison table
#    whose style switches from the deck's custom "Table_0" style to the
#    built-in table style {689E48D7-6EA6-481B-B047-EF69FE79BBDC}.
# 2) The deck's theme switches from the custom "Integral" theme to the
#    stock "Office Theme" that was already embedded (until now only used
#    by the notes master) - i.e. the Design applied to the slide master
#    changes.
# -----------------------------------------------------------------------
$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 ---
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{689E48D7-6EA6-481B-B047-EF69FE79BBDC}")

# --- 2) Swap the active deck theme ("Integral" -> "Office Theme") ---
# The "Office Theme" theme part is already in the package (it backs the
# notes master); applying it to the slide master makes it the deck's
# active theme, while the previous "Integral" theme becomes the one kept
# around for the notes master.
$slideMaster = $p.Slides.Item(1).Master
$notesMaster = $p.NotesMaster
$slideMaster.ApplyTheme("ppt/theme/theme2.xml")
$notesMaster.ApplyTheme("ppt/theme/theme1.xml")
